# Relation Schemas.pptx — "Add files via upload" commit
#
# The only change that applies to this single-slide excerpt is in the
# bulleted "TextBox 3" shape on slide 1: the relation
#     ordAdd(order_ num, add_id)
# gained two more (underlined) attributes, becoming
#     ordAdd(order_ num, add_id, isShipping, isBilling)
#
# (The diff also touched the per-slide "datetimeFigureOut" footer field on
#  many other slides, bumping it from 2020-04-01 to 2020-04-03, but none of
#  those slides are part of this single-slide deck, so there is nothing to
#  change there.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape that holds the bulleted relation-schema list (it has the
# "ordAdd(...)" bullet among its paragraphs) instead of hard-coding a shape
# index, so the script is resilient to shape ordering.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*ordAdd(order_ num, add_id)*") {
            $targetShape = $candidate
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the shape containing 'ordAdd(order_ num, add_id)'"
}

$tr = $targetShape.TextFrame.TextRange

# Locate the exact paragraph "* ordAdd(order_ num, add_id)" among the
# paragraphs of the text box.
$paraCount = $tr.Paragraphs().Count
$paraIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "*ordAdd(order_ num, add_id)*") {
        $paraIndex = $i
    }
}

if ($paraIndex -eq -1) {
    throw "Could not find paragraph 'ordAdd(order_ num, add_id)'"
}

$para = $tr.Paragraphs($paraIndex, 1)
$paraStart = $para.Start
$paraTextBefore = $para.Text

# Offset (0-based, relative to paragraph start) of the closing ")" that
# follows "add_id".
$closeParenOffset = $paraTextBefore.IndexOf("add_id)") + ("add_id".Length)

$closeParen = $tr.Characters($paraStart + $closeParenOffset, 1)

# Extend "...add_id)" into "...add_id, isShipping, isBilling)". The newly
# typed text inherits the non-underlined italic formatting of the ")" run
# it replaces; underlining of the individual new identifiers is applied
# below, matching the source formatting of "order_ num" / "add_id".
$closeParen.Text = ", isShipping, isBilling)"

# Re-fetch the paragraph now that its length has changed.
$para2 = $tr.Paragraphs($paraIndex, 1)
$start2 = $para2.Start

# The new suffix starts exactly where the old ")" used to be (the ")" run
# got overwritten in place), so reuse that already-known offset instead of
# re-deriving it from the (carriage-return padded) .Text length.
$suffixOffset = $closeParenOffset

# ", " before isShipping
$offComma1 = $suffixOffset
$r1 = $tr.Characters($start2 + $offComma1, 2)
$r1.Font.Underline = -1

# "isShipping"
$offShipping = $suffixOffset + 2
$r2 = $tr.Characters($start2 + $offShipping, 10)
$r2.Font.Underline = -1

# ", " before isBilling
$offComma2 = $offShipping + 10
$r3 = $tr.Characters($start2 + $offComma2, 2)
$r3.Font.Underline = -1

# "isBilling"
$offBilling = $offComma2 + 2
$r4 = $tr.Characters($start2 + $offBilling, 9)
$r4.Font.Underline = -1
